$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.399.76"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").Value = "1.736.85"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.37"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.34"
$ws.Range("E8").Value = "  +13.22%  "
$ws.Range("E9").Value = "  +4.35%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "1.977.40"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "1.735.94"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.85"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "28.359.89"
$ws.Range("E17").Value = "  +4.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.67"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.04"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.73"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.42"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").Value = "1.504.93"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.608"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  +1.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.08"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.93"
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "1.882.05"
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.806"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("E47").Value = "  +9.73%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "91.20"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0112"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.27"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("E51").Value = "  +1.48%  "
